$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 535.2
$ws.Range("I9").Value = 535.2
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 535.2
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -366.2
$ws.Range("N9").ClearContents()

$ws.Range("H12").Value = 747.5
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 896.6667
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 896.6667
$ws.Range("M12").Value = -130
$ws.Range("N12").Value = -1236.6667

$ws.Range("H38").Value = 550
$ws.Range("I38").Value = 433.33334
$ws.Range("K38").Value = 1300.00002
$ws.Range("M38").Value = -928.0000199999999

$ws.Range("H58").Value = 1441.4286
$ws.Range("I58").Value = 1018
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 3054
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -2904
$ws.Range("N58").Value = -7800

$ws.Range("H87").Value = 24428.285
$ws.Range("J87").Value = 24428.285
$ws.Range("L87").Value = 24428.285
$ws.Range("N87").Value = -26924.285

$ws.Range("H90").Value = 24428.285
$ws.Range("J90").Value = 24428.285
$ws.Range("L90").Value = 73284.855
$ws.Range("N90").Value = -85764.855

$ws.Range("H112").Value = 5471.4707
$ws.Range("I112").Value = 366.66666
$ws.Range("J112").Value = 5965.484
$ws.Range("K112").Value = 1099.99998
$ws.Range("L112").Value = 17896.452
$ws.Range("M112").Value = 8.00001999999995
$ws.Range("N112").Value = -20112.452

$ws.Range("H125").Value = 1257.2858
$ws.Range("I125").Value = 446.2
$ws.Range("J125").Value = 1707.8889
$ws.Range("K125").Value = 4015.8
$ws.Range("L125").Value = 15371.0001
$ws.Range("M125").Value = -1555.8
$ws.Range("N125").Value = -20291.0001

$ws.Range("H138").Value = 1319.5
$ws.Range("I138").Value = 544.55554
$ws.Range("J138").Value = 1953.5454
$ws.Range("K138").Value = 1633.66662
$ws.Range("L138").Value = 5860.6362
$ws.Range("M138").Value = 3506.33338
$ws.Range("N138").Value = -16140.6362

$ws.Range("H141").Value = 2483.1904
$ws.Range("I141").Value = 818.1458
$ws.Range("J141").Value = 7811.3335
$ws.Range("K141").Value = 2454.4374
$ws.Range("L141").Value = 23434.0005
$ws.Range("M141").Value = 2725.5626
$ws.Range("N141").Value = -33794.00049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1787.7858
$ws.Range("I2").Value = 1811.25
$ws.Range("J2").Value = 1756.5
$ws.Range("K2").Value = 1811.25
$ws.Range("L2").Value = 1756.5
$ws.Range("M2").Value = -1698.25
$ws.Range("N2").Value = -1982.5

$ws.Range("H32").Value = 797839.5600000001
$ws.Range("I32").Value = 924251.5
$ws.Range("J32").Value = 16384.273
$ws.Range("K32").Value = 924251.5
$ws.Range("L32").Value = 16384.273
$ws.Range("M32").Value = -923964.5
$ws.Range("N32").Value = -16958.273

$ws.Range("H45").Value = 1806.3043
$ws.Range("I45").Value = 1791.875
$ws.Range("J45").Value = 1839.2858
$ws.Range("K45").Value = 1791.875
$ws.Range("L45").Value = 1839.2858
$ws.Range("M45").Value = -1414.875
$ws.Range("N45").Value = -2593.2858

$ws.Range("H61").Value = 1634.2554
$ws.Range("I61").Value = 1320.3914
$ws.Range("J61").Value = 2500.52
$ws.Range("K61").Value = 1320.3914
$ws.Range("L61").Value = 2500.52
$ws.Range("M61").Value = -1108.3914
$ws.Range("N61").Value = -2924.52

$ws.Range("H75").Value = 70000
$ws.Range("J75").Value = 70000
$ws.Range("L75").Value = 70000
$ws.Range("N75").Value = -71748

$ws.Range("H78").Value = 70000
$ws.Range("J78").Value = 70000
$ws.Range("L78").Value = 210000
$ws.Range("N78").Value = -218736

$ws.Range("H116").Value = 1787.7858
$ws.Range("I116").Value = 1811.25
$ws.Range("J116").Value = 1756.5
$ws.Range("K116").Value = 1811.25
$ws.Range("L116").Value = 1756.5
$ws.Range("M116").Value = 482.75
$ws.Range("N116").Value = -6344.5

$ws.Range("H122").Value = 92006.73
$ws.Range("I122").Value = 111895.555
$ws.Range("K122").Value = 335686.665
$ws.Range("M122").Value = -333236.665

$ws.Range("H136").Value = 1634.2554
$ws.Range("I136").Value = 1320.3914
$ws.Range("J136").Value = 2500.52
$ws.Range("K136").Value = 3961.1742
$ws.Range("L136").Value = 7501.559999999999
$ws.Range("M136").Value = -1411.1742
$ws.Range("N136").Value = -12601.56

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1787.7858
$ws.Range("I3").Value = 1811.25
$ws.Range("J3").Value = 1756.5
$ws.Range("K3").Value = 1811.25
$ws.Range("L3").Value = 1756.5
$ws.Range("M3").Value = -1697.25
$ws.Range("N3").Value = -1984.5

$ws.Range("H107").Value = 78273.69500000001
$ws.Range("I107").Value = 112116.11
$ws.Range("J107").Value = 2128.25
$ws.Range("K107").Value = 112116.11
$ws.Range("L107").Value = 2128.25
$ws.Range("M107").Value = -110196.11
$ws.Range("N107").Value = -5968.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 4167615.8
$ws.Range("I107").Value = 10417240
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 10417240
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = -10415320
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1112.2667
$ws.Range("J5").Value = 1674.7693
$ws.Range("L5").Value = 5024.3079
$ws.Range("N5").Value = -5248.3079

$ws.Range("H64").Value = 1900.6
$ws.Range("I64").Value = 1337.3334
$ws.Range("K64").Value = 4012.0002
$ws.Range("M64").Value = -3742.0002

$ws.Range("H67").Value = 1900.6
$ws.Range("I67").Value = 1337.3334
$ws.Range("K67").Value = 4012.0002
$ws.Range("M67").Value = -3076.0002

$ws.Range("H92").Value = 605.5714
$ws.Range("I92").Value = 396.33334
$ws.Range("J92").Value = 762.5
$ws.Range("K92").Value = 1189.00002
$ws.Range("L92").Value = 2287.5
$ws.Range("M92").Value = 58.99998000000005
$ws.Range("N92").Value = -4783.5

$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242

$ws.Range("H122").Value = 3383.25
$ws.Range("I122").Value = 413.72
$ws.Range("J122").Value = 10132.182
$ws.Range("K122").Value = 3723.48
$ws.Range("L122").Value = 91189.63800000001
$ws.Range("M122").Value = -1273.48
$ws.Range("N122").Value = -96089.63800000001

$ws.Range("H131").Value = 2490.2666
$ws.Range("I131").Value = 532.6667
$ws.Range("J131").Value = 2660.4927
$ws.Range("K131").Value = 1598.0001
$ws.Range("L131").Value = 7981.478099999999
$ws.Range("M131").Value = 3441.9999
$ws.Range("N131").Value = -18061.4781

$ws.Range("H135").Value = 1112.2667
$ws.Range("J135").Value = 1674.7693
$ws.Range("L135").Value = 15072.9237
$ws.Range("N135").Value = -20142.9237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1675.6
$ws.Range("I102").Value = 1662.0834
$ws.Range("K102").Value = 1662.0834
$ws.Range("M102").Value = -40.08339999999998

$ws.Range("H126").Value = 1881.5
$ws.Range("I126").Value = 1806
$ws.Range("J126").Value = 1957
$ws.Range("K126").Value = 5418
$ws.Range("L126").Value = 5871
$ws.Range("M126").Value = -2948
$ws.Range("N126").Value = -10811

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 95000
$ws.Range("J76").Value = 95000
$ws.Range("L76").Value = 95000
$ws.Range("N76").Value = -95676

$ws.Range("H79").Value = 95000
$ws.Range("J79").Value = 95000
$ws.Range("L79").Value = 95000
$ws.Range("N79").Value = -97340

$ws.Range("H136").Value = 3877376.2
$ws.Range("I136").Value = 1363.9412
$ws.Range("J136").Value = 18520090
$ws.Range("K136").Value = 4091.8236
$ws.Range("L136").Value = 55560270
$ws.Range("M136").Value = -1541.8236
$ws.Range("N136").Value = -55565370

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9823.25
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 9823.25
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 9823.25
$ws.Range("N54").Value = -10863.25
$ws.Range("M54").ClearContents()

$ws.Range("H81").Value = 3153.5
$ws.Range("I81").Value = 3555.8333
$ws.Range("J81").Value = 2550
$ws.Range("K81").Value = 7111.6666
$ws.Range("L81").Value = 5100
$ws.Range("M81").Value = -6050.6666
$ws.Range("N81").Value = -7222

$ws.Range("H84").Value = 3153.5
$ws.Range("I84").Value = 3555.8333
$ws.Range("J84").Value = 2550
$ws.Range("K84").Value = 35558.333
$ws.Range("L84").Value = 25500
$ws.Range("M84").Value = -30254.333
$ws.Range("N84").Value = -36108

$ws.Range("H107").Value = 632.91895
$ws.Range("I107").Value = 612.9666999999999
$ws.Range("J107").Value = 718.4286
$ws.Range("K107").Value = 1838.9001
$ws.Range("L107").Value = 2155.2858
$ws.Range("M107").Value = 81.09990000000016
$ws.Range("N107").Value = -5995.2858
